# Insert two new weekly price rows (Conconina(o) / Escarola) for the
# "Terminal Hortofrutícola Agro Chillán - Lechuga" sheet.
#
# The new rows are inserted right before the existing row 493, pushing
# the former rows 493-523 down to 495-525 (dimension grows from R523 to
# R525). The two new rows carry the same "shape" as the old 493/494 rows
# (same market/region/category columns) but with an updated date and a
# few updated price/volume figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 493 (shifts old 493.. down by 2, carries the
# date-style formatting on column D along, just like Excel's own
# "Insert Copied/Above Cells" behaviour).
$ws.Range("A493:R494").EntireRow.Insert()

# New row 493: Conconina(o), 2022-01-24 (serial 44585)
$ws.Cells.Item(493, 1).Value = 7
$ws.Cells.Item(493, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(493, 3).Value = "Ñuble"
$ws.Cells.Item(493, 4).Value = 44585
$ws.Cells.Item(493, 5).Value = 16
$ws.Cells.Item(493, 6).Value = 100112033
$ws.Cells.Item(493, 7).Value = "Lechuga"
$ws.Cells.Item(493, 8).Value = "Conconina(o)"
$ws.Cells.Item(493, 9).Value = "Primera"
$ws.Cells.Item(493, 10).Value = 120
$ws.Cells.Item(493, 11).Value = 5000
$ws.Cells.Item(493, 12).Value = 5500
$ws.Cells.Item(493, 13).Value = 5250
$ws.Cells.Item(493, 14).Value = "$/caja 10 unidades"
$ws.Cells.Item(493, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(493, 16).Value = 525
$ws.Cells.Item(493, 17).Value = 10
$ws.Cells.Item(493, 18).Value = "Hortaliza"

# New row 494: Escarola, 2022-01-24 (serial 44585)
$ws.Cells.Item(494, 1).Value = 7
$ws.Cells.Item(494, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(494, 3).Value = "Ñuble"
$ws.Cells.Item(494, 4).Value = 44585
$ws.Cells.Item(494, 5).Value = 16
$ws.Cells.Item(494, 6).Value = 100112033
$ws.Cells.Item(494, 7).Value = "Lechuga"
$ws.Cells.Item(494, 8).Value = "Escarola"
$ws.Cells.Item(494, 9).Value = "Primera"
$ws.Cells.Item(494, 10).Value = 120
$ws.Cells.Item(494, 11).Value = 6000
$ws.Cells.Item(494, 12).Value = 6500
$ws.Cells.Item(494, 13).Value = 6250
$ws.Cells.Item(494, 14).Value = "$/caja 15 unidades"
$ws.Cells.Item(494, 15).Value = "Región del Maule"
$ws.Cells.Item(494, 16).Value = 417
$ws.Cells.Item(494, 17).Value = 15
$ws.Cells.Item(494, 18).Value = "Hortaliza"
